$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 is TestScenario_3 / "Edit Account" test case.
# The reviewer changed the Approved/Rejected verdict from "Approved" to
# "Rejected" and filled in a reason in the ReasonToReject column.
$ws.Range("I8").Value = "Rejected"
$ws.Range("J8").Value = "Not required"

# Reflect the resulting scroll position / selection from the edit session.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J14").Select()
